$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 89585255
$ws.Range("B2").Value = 77506
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = 'Garnlav'
$ws.Range("G2").Value = 'Alectoria sarmentosa'
$ws.Range("H2").Value = '(Ach.) Ach.'
$ws.Range("Q2").Value = 422268.873820144
$ws.Range("R2").Value = 6941277.034413425

$ws.Range("A3").Value = 89585290
$ws.Range("B3").Value = 73693
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6440
$ws.Range("F3").Value = 'Vitgrynig nållav'
$ws.Range("G3").Value = 'Chaenotheca subroscida'
$ws.Range("H3").Value = '(Eitner) Zahlbr.'
$ws.Range("Q3").Value = 422382.1566529368
$ws.Range("R3").Value = 6941443.959308082

$ws.Range("A4").Value = 89585224
$ws.Range("B4").Value = 89356
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = 'Vedticka'
$ws.Range("G4").Value = 'Fuscoporia viticola'
$ws.Range("H4").Value = '(Schwein.) Murrill'
$ws.Range("Q4").Value = 422227.8240813478
$ws.Range("R4").Value = 6941277.997683762

$ws.Range("A5").Value = 89585286
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = 'Garnlav'
$ws.Range("G5").Value = 'Alectoria sarmentosa'
$ws.Range("H5").Value = '(Ach.) Ach.'
$ws.Range("Q5").Value = 422195.0116517335
$ws.Range("R5").Value = 6941276.003047519

$ws.Range("A6").Value = 89585289
$ws.Range("B6").Value = 77588
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 864
$ws.Range("F6").Value = 'Knottrig blåslav'
$ws.Range("G6").Value = 'Hypogymnia bitteri'
$ws.Range("H6").Value = '(Lynge) Ahti'
$ws.Range("Q6").Value = 422249.1704919123
$ws.Range("R6").Value = 6941283.026650068

$ws.Range("A7").Value = 89585233
$ws.Range("B7").Value = 77506
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'
$ws.Range("Q7").Value = 422310.9247597085
$ws.Range("R7").Value = 6941299.089533247

$ws.Range("A8").Value = 89585234
$ws.Range("B8").Value = 77506
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = 'Garnlav'
$ws.Range("G8").Value = 'Alectoria sarmentosa'
$ws.Range("H8").Value = '(Ach.) Ach.'
$ws.Range("Q8").Value = 422440.1516844163
$ws.Range("R8").Value = 6941536.143204342

$ws.Range("A9").Value = 89585240
$ws.Range("B9").Value = 77588
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 864
$ws.Range("F9").Value = 'Knottrig blåslav'
$ws.Range("G9").Value = 'Hypogymnia bitteri'
$ws.Range("H9").Value = '(Lynge) Ahti'
$ws.Range("Q9").Value = 422342.0538613191
$ws.Range("R9").Value = 6941308.037107519

$ws.Range("A10").Value = 89585283
$ws.Range("B10").Value = 78596
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 6462
$ws.Range("F10").Value = 'Stuplav'
$ws.Range("G10").Value = 'Nephroma bellum'
$ws.Range("H10").Value = '(Spreng.) Tuck.'
$ws.Range("Q10").Value = 422291.9279179227
$ws.Range("R10").Value = 6941295.848404294

$ws.Range("A11").Value = 89585261
$ws.Range("B11").Value = 77588
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 864
$ws.Range("F11").Value = 'Knottrig blåslav'
$ws.Range("G11").Value = 'Hypogymnia bitteri'
$ws.Range("H11").Value = '(Lynge) Ahti'
$ws.Range("Q11").Value = 422513.8717266123
$ws.Range("R11").Value = 6941570.820536849
